# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The worksheet's monthly-arrears table (column E, rows 16-24) is refreshed:
# the nine "Periodo Mora" values are re-entered in descending (most-recent-first)
# order instead of the previous ascending order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periodos = @("2306", "2305", "2304", "2303", "2302", "2301", "2212", "2211", "2210")

$row = 16
foreach ($periodo in $periodos) {
    $ws.Range("E" + $row).Value = $periodo
    $row++
}
